# Moved new classifier into java resources folder. Added per dataset training results.
# This script adds two new columns (H, I) of "Individual training performance" results
# to the "Foaie2" worksheet, with appropriate headers, merges, row heights and widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Activate()

$NL = [char]10

# ---------------------------------------------------------------------------
# Column widths for the two new columns (H, I)
# ---------------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 29.8
$ws.Columns.Item(9).ColumnWidth = 29.8

# ---------------------------------------------------------------------------
# Row height tweaks (row 5 gains an explicit custom height flag, row 8 shrinks
# from the old 45.75 auto height down to a fixed 30)
# ---------------------------------------------------------------------------
$ws.Rows.Item(5).RowHeight = 30
$ws.Rows.Item(8).RowHeight = 30

# ---------------------------------------------------------------------------
# Headers (row 1/2, merged like the existing G1:G2 header)
# ---------------------------------------------------------------------------
$ws.Range("H1:H2").Merge()
$ws.Range("H1").Value = "Individual training performance (F1 at 20 epochs)"

$ws.Range("I1:I2").Merge()
$ws.Range("I1").Value = "Individual training performance on orig DS (F1 at 20 epochs)"

foreach ($addr in @("H1", "I1")) {
    $c = $ws.Range($addr)
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4108
    $c.WrapText = $true
}

# ---------------------------------------------------------------------------
# Data rows 3-15 (one value per row in H; I only populated on a few rows,
# with I9:I15 merged into a single "N/A" cell)
# ---------------------------------------------------------------------------
$hValues = @{
    3  = "0.89 0.89 0.89${NL}0.81 0.82 0.82"
    4  = "0.90 0.97 0.94${NL}0.80 0.55 0.65"
    5  = "0.82 0.84 0.83${NL}0.64 0.60 0.62"
    6  = "0.84 0.93 0.88${NL}0.76 0.55 0.64"
    7  = "0.84 0.94 0.89${NL}0.77 0.53 0.63"
    8  = "0.85 0.93 0.89${NL}0.76 0.56 0.65"
    9  = "0.84 0.92 0.88${NL}0.69 0.50 0.58"
    10 = "0.80 0.94 0.86${NL}0.65 0.25 0.36"
    11 = "0.83 0.90 0.86${NL}0.74 0.59 0.66"
    12 = "0.85 0.90 0.87${NL}0.76 0.67 0.71"
    13 = "0.77 0.92 0.84${NL}0.76 0.47 0.58"
    14 = "0.82 0.90 0.86${NL}0.68 0.50 0.58"
    15 = "0.85 0.92 0.88${NL}0.62 0.43 0.50"
}

foreach ($row in $hValues.Keys) {
    $ws.Range("H$row").Value = $hValues[$row]
}

$ws.Range("I4").Value = "0.86 0.97 0.91${NL}0.83 0.49 0.62"
$ws.Range("I5").Value = "0.90 0.97 0.93${NL}0.70 0.40 0.51"

$ws.Range("I9:I15").Merge()
$ws.Range("I9").Value = "N/A"

# ---------------------------------------------------------------------------
# Rows 16-21 -- both H and I collapse to a single merged "N/A" cell, mirroring
# the existing G16:G21 merge
# ---------------------------------------------------------------------------
$ws.Range("H16:H21").Merge()
$ws.Range("H16").Value = "N/A"

$ws.Range("I16:I21").Merge()
$ws.Range("I16").Value = "N/A"

# ---------------------------------------------------------------------------
# Alignment / wrap text for all the new data cells
# ---------------------------------------------------------------------------
foreach ($addr in @("H3:H15", "I3:I8", "I9:I15", "H16:H21", "I16:I21")) {
    $c = $ws.Range($addr)
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4108
    $c.WrapText = $true
}

# ---------------------------------------------------------------------------
# Borders: a medium line separates the new block from the rest of the table
# (left edge of H), a medium line closes the table off on the right (right
# edge of I), and thin lines separate the individual per-row entries, mirroring
# the existing G-column box look.
# ---------------------------------------------------------------------------
$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10
$xlContinuous = 1
$xlThin = 2
$xlMedium = -4138

$full = $ws.Range("H1:I21")
$full.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
$full.Borders.Item($xlEdgeTop).Weight = $xlMedium
$full.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
$full.Borders.Item($xlEdgeBottom).Weight = $xlMedium

$hCol = $ws.Range("H1:H21")
$hCol.Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
$hCol.Borders.Item($xlEdgeLeft).Weight = $xlMedium

$iCol = $ws.Range("I1:I21")
$iCol.Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
$iCol.Borders.Item($xlEdgeRight).Weight = $xlMedium

foreach ($addr in @("H1:I2", "H3:I8", "H9:I15", "H16:I21")) {
    $blk = $ws.Range($addr)
    $blk.Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
    $blk.Borders.Item($xlEdgeTop).Weight = $xlMedium
    $blk.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
    $blk.Borders.Item($xlEdgeBottom).Weight = $xlMedium
}

foreach ($row in 3..14) {
    $r = $ws.Range("H$row`:I$row")
    $r.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
    $r.Borders.Item($xlEdgeBottom).Weight = $xlThin
}

# ---------------------------------------------------------------------------
# Sheet view: drop the old scroll position / selection, zoom to 85% and select H15
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 85
$ws.Range("H15").Select()
